$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force all touched cells to Text format first so Excel does not
# auto-convert numeric-looking strings (e.g. "1.000", "0.9990") into numbers,
# matching the original inlineStr/text semantics of the sheet.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.383.48'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.849.49'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.50'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6285'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07589'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2919'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.54'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07751'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.849.73'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.014'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6787'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001043'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.20'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.101.14'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.110'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.379.33'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '230.13'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.72%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.35'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.431'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9999'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.33'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.12%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1395'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.20%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.443'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.68'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.426'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +6.01%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.474'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05685'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.119'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.047'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.156'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.68%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.824'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7004'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.583'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01825'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.78%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.237.86'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.45%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.718'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.421'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.96%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9030'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9997'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.008.14'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.23%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.47'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.36%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.77'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.139'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1157'
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3986'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.999'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.681'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.05%  '
